# Add four newly collected submissions to the first worksheet
# ("八位序列号收集收集结果yd5"), continuing the existing log of
# submitter name / submit time / collected ID rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFormat = "yyyy/m/d h:mm:ss;@"

$newRows = @(
    @{ Row = 21; Name = "知道分子.";                      Time = 45897.007962963;  Id = "2b5431f5" },
    @{ Row = 22; Name = "MEING";                          Time = 45897.4482060185; Id = "44d0e155" },
    @{ Row = 23; Name = "-";                               Time = 45897.7506134259; Id = "a04f3a54" },
    @{ Row = 24; Name = "　　　　　　　　　　　　";        Time = 45897.9142013889; Id = "590320c9" }
)

foreach ($entry in $newRows) {
    $r = $entry.Row
    $ws.Range("A$r").Value = $entry.Name
    $ws.Range("B$r").Value = $entry.Time
    $ws.Range("B$r").NumberFormat = $dateFormat
    $ws.Range("C$r").Value = $entry.Id
}
